# JS + html flow from MADtrial
# Adds an "if / begin screen / note / end screen / end if" block (x2) to the
# "survey" sheet, replacing the old placeholder "This is an example" note
# with a proper info screen, and driving the two new screens from the
# ERROR1 / ERROR2 data values.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Replace the old placeholder note text with the new info screen text.
$survey.Range("G3").Value = "Info screen"
$survey.Range("H3").Value = "Info screen"

# --- First error-handling block skeleton (rows 5-11) ---
$survey.Range("B5").Value = "if"
$survey.Range("C5").Value = 'data("ERROR1") == "1"'

$survey.Range("B6").Value = "begin screen"
$survey.Range("B10").Value = "end screen"
$survey.Range("B11").Value = "end if"

# --- Second error-handling block skeleton (rows 12-18) ---
$survey.Range("B12").Value = "if"
$survey.Range("C12").Value = 'data("ERROR2") == "2"'

$survey.Range("B13").Value = "begin screen"
$survey.Range("B17").Value = "end screen"
$survey.Range("B18").Value = "end if"

# --- Fill in the note contents for each screen ---
$survey.Range("D7").Value = "note"
$survey.Range("G7").Value = "Cleaning of error 1"
$survey.Range("H7").Value = "Cleaning of error 1"

$survey.Range("D14").Value = "note"
$survey.Range("G14").Value = "Cleaning of error 2"
$survey.Range("H14").Value = "Cleaning of error 2"

# Move the active selection to reflect where editing left off.
$choices = $wb.Worksheets.Item("choices")
$choices.Range("B21").Select()

$survey.Range("D15").Select()
